# Apply the "regular linear regression" update: new resampled IDs for the
# "Aerial Grounding" sheet, refreshed hazard flags, corrected SUM formula,
# and an updated Summary sheet (new accuracy numbers + a new STDEV.S column).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Aerial Grounding" sheet: replace the 15 sampled IDs (A2:A16) with a
#    newly drawn sample, and refresh the "Contains Hazard" flags (B2:B16).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Aerial Grounding")

$newIds = @(
    "2010_AK-MSS-001278_EKLUTNA_0",
    "2008_AZ-PPA-000056_SOLANO_0",
    "2012_WA-OWF-000610_OKANOGAN COMPLEX_6",
    "2013_CA-KNF-005949_FORKS COMPLEX_67",
    "2013_MT-SWS-000068_LOLO CREEK COMPLEX_4",
    "2011_AZ-CNF-011047_HORSESHOE 2_30",
    "2012_WA-OWF-000610_OKANOGAN COMPLEX_5",
    "2006_WY-SHF-000152_BOMBER BASIN_0",
    "2013_MT-SWS-000068_LOLO CREEK COMPLEX_3",
    "2006_CA-SHF-001693_BAR COMPLEX_110",
    "2006_MT-GNF-055_PARADISE VALLEY COMPLEX_50",
    "2013_ID-SCF-13165_LODGEPOLE_7",
    "2006_WY-SHF-000152_BOMBER BASIN_1",
    "2012_UT-NWS-000507_FLOOD CANYON_4",
    "2013_CA-STF-002857_RIM_15"
)

$newFlags = @(1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 0, 1, 0)

for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value2 = $newIds[$i]
    $ws.Range("B$row").Value2 = $newFlags[$i]
}

# Re-enter the SUM formula (now without the stray leading space) so it
# recalculates to the new total (12).
$ws.Range("B17").Formula = "=SUM(B2:B16)"

# Update the sheet's selection to match the saved view state.
$ws.Activate()
$ws.Range("A1:B18").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. "Summary" sheet: refresh the Aerial Grounding row's correct-sample
#    count, add the new standard-deviation column, and update the
#    selection to match the saved view state.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("E17").Value2 = 12

$summary.Range("G20").Formula = "=STDEV.S(F2:F18)"

$summary.Activate()
$summary.Range("G20").Select() | Out-Null
